$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2..10 (A..T), per Dr Hou's advice: the
# Tnfsf13->Tnfrsf11b sending/target cluster cross (ECs / FAPs / sCs) is
# recomputed and expanded from a 2x2-ish partial table into the full 3x3
# cross of clusters.
$data = @(
    @("ECs","Tnfsf13","Tnfrsf11b","ECs",3,1,0.291812,0.875436,0.1179458223311005,0.1179458223311005,1,0.3333333333333333,0.1176943333333333,0.353083,0.05829606481791055,0.05829606481791055,0.03434461879866667,0.309101569188,0.006875777303615596,0.006875777303615596),
    @("ECs","Tnfsf13","Tnfrsf11b","FAPs",3,1,0.291812,0.875436,0.1179458223311005,0.1179458223311005,3,1,1.183046666666667,3.54914,0.5859837360842607,0.5859837360842608,0.3452272138933333,3.10704492504,0.0691143336251087,0.0691143336251087),
    @("ECs","Tnfsf13","Tnfrsf11b","sCs",3,1,0.291812,0.875436,0.1179458223311005,0.1179458223311005,3,1,0.718166,2.154498,0.3557201990978286,0.3557201990978286,0.209569456792,1.886125111128,0.04195571140237619,0.04195571140237619),
    @("FAPs","Tnfsf13","Tnfrsf11b","ECs",2,0.6666666666666666,0.514066,1.542198,0.2077773947009016,0.2077773947009016,1,0.3333333333333333,0.1176943333333333,0.353083,0.05829606481791055,0.05829606481791055,0.06050265515933333,0.544523896434,0.01211260446918034,0.01211260446918034),
    @("FAPs","Tnfsf13","Tnfrsf11b","FAPs",2,0.6666666666666666,0.514066,1.542198,0.2077773947009016,0.2077773947009016,3,1,1.183046666666667,3.54914,0.5859837360842607,0.5859837360842608,0.6081640677466666,5.47347660972,0.1217541740206884,0.1217541740206884),
    @("FAPs","Tnfsf13","Tnfrsf11b","sCs",2,0.6666666666666666,0.514066,1.542198,0.2077773947009016,0.2077773947009016,3,1,0.718166,2.154498,0.3557201990978286,0.3557201990978286,0.369184722956,3.322662506603999,0.07391061621103286,0.07391061621103286),
    @("sCs","Tnfsf13","Tnfrsf11b","ECs",3,1,1.668241,5.004723,0.674276782967998,0.674276782967998,1,0.3333333333333333,0.1176943333333333,0.353083,0.05829606481791055,0.05829606481791055,0.1963425123343333,1.767082611009,0.03930768304511461,0.03930768304511462),
    @("sCs","Tnfsf13","Tnfrsf11b","FAPs",3,1,1.668241,5.004723,0.674276782967998,0.674276782967998,3,1,1.183046666666667,3.54914,0.5859837360842607,0.5859837360842608,1.973606954246667,17.76246258822,0.3951152284384636,0.3951152284384638),
    @("sCs","Tnfsf13","Tnfrsf11b","sCs",3,1,1.668241,5.004723,0.674276782967998,0.674276782967998,3,1,0.718166,2.154498,0.3557201990978286,0.3557201990978286,1.198073966006,10.782665694054,0.2398538714844196,0.2398538714844196)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
